# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# Add the new Wild Card round result row for C.Saunders on the WR sheet
$wsWR = $wb.Worksheets.Item("WR")
$wsWR.Range("A10").Value = "C.Saunders"
$wsWR.Range("B10:J10").Value = 0

# Update the cursor/selection left on the RB sheet after the edit
$wsRB = $wb.Worksheets.Item("RB")
$wsRB.Range("A7").Select()

# Move the cursor/selection on the WR sheet to below the newly logged row,
# and leave WR as the active (selected) tab
$wsWR.Range("J11").Select()
$wsWR.Activate()
